$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 34148744
$ws.Range("I100").Value = 41680676
$ws.Range("J100").Value = 22850846
$ws.Range("K100").Value = 41680676
$ws.Range("L100").Value = 22850846
$ws.Range("M100").Value = -41680135
$ws.Range("N100").Value = -22851928
$ws.Range("H107").Value = 6738.1035
$ws.Range("I107").Value = 6939.1177
$ws.Range("J107").Value = 6453.3335
$ws.Range("K107").Value = 6939.1177
$ws.Range("L107").Value = 6453.3335
$ws.Range("M107").Value = -5019.1177
$ws.Range("N107").Value = -10293.3335
$ws.Range("H111").Value = 6299.5
$ws.Range("J111").Value = 7000
$ws.Range("L111").Value = 21000
$ws.Range("N111").Value = -27134
$ws.Range("H112").Value = 51570.227
$ws.Range("J112").Value = 80299.64
$ws.Range("L112").Value = 240898.92
$ws.Range("N112").Value = -243114.92
$ws.Range("H116").Value = 623580.5600000001
$ws.Range("I116").Value = 1394122.5
$ws.Range("K116").Value = 1394122.5
$ws.Range("M116").Value = -1390680.5
$ws.Range("H125").Value = 3553
$ws.Range("I125").Value = 2141.3333
$ws.Range("K125").Value = 19271.9997
$ws.Range("M125").Value = -16811.9997
$ws.Range("H127").Value = 2580.2144
$ws.Range("I127").Value = 1109.4615
$ws.Range("K127").Value = 3328.3845
$ws.Range("M127").Value = 1631.6155
$ws.Range("H135").Value = 3767.1333
$ws.Range("I135").Value = 3767.1333
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 33904.1997
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -31369.1997
$ws.Range("N135").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13002.4
$ws.Range("I2").Value = 16467.715
$ws.Range("K2").Value = 16467.715
$ws.Range("M2").Value = -16354.715
$ws.Range("H45").Value = 13476.125
$ws.Range("J45").Value = 5831.9
$ws.Range("L45").Value = 5831.9
$ws.Range("N45").Value = -6585.9
$ws.Range("H61").Value = 7759.25
$ws.Range("I61").Value = 12121
$ws.Range("K61").Value = 12121
$ws.Range("M61").Value = -11909
$ws.Range("H74").Value = 3079.853
$ws.Range("I74").Value = 1552.3214
$ws.Range("J74").Value = 10208.333
$ws.Range("K74").Value = 1552.3214
$ws.Range("L74").Value = 10208.333
$ws.Range("M74").Value = -678.3214
$ws.Range("N74").Value = -11956.333
$ws.Range("H77").Value = 3079.853
$ws.Range("I77").Value = 1552.3214
$ws.Range("J77").Value = 10208.333
$ws.Range("K77").Value = 7761.607
$ws.Range("L77").Value = 51041.665
$ws.Range("M77").Value = -3393.607
$ws.Range("N77").Value = -59777.665
$ws.Range("H110").Value = 1935.591
$ws.Range("I110").Value = 1710.9412
$ws.Range("K110").Value = 1710.9412
$ws.Range("M110").Value = 334.0588
$ws.Range("H116").Value = 13002.4
$ws.Range("I116").Value = 16467.715
$ws.Range("K116").Value = 16467.715
$ws.Range("M116").Value = -14173.715
$ws.Range("H132").Value = 3264.125
$ws.Range("I132").Value = 3014.9443
$ws.Range("K132").Value = 9044.832900000001
$ws.Range("M132").Value = -6514.832900000001
$ws.Range("H136").Value = 7759.25
$ws.Range("I136").Value = 12121
$ws.Range("K136").Value = 36363
$ws.Range("M136").Value = -33813
$ws.Range("H138").Value = 65001.89
$ws.Range("J138").Value = 65001.89
$ws.Range("L138").Value = 65001.89
$ws.Range("N138").Value = -75281.89

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13002.4
$ws.Range("I3").Value = 16467.715
$ws.Range("K3").Value = 16467.715
$ws.Range("M3").Value = -16353.715
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51622
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -158112
$ws.Range("H134").Value = 2927.7693
$ws.Range("I134").Value = 2868.5454
$ws.Range("J134").Value = 3253.5
$ws.Range("K134").Value = 8605.636200000001
$ws.Range("L134").Value = 9760.5
$ws.Range("M134").Value = -6070.636200000001
$ws.Range("N134").Value = -14830.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 110.36364
$ws.Range("I7").Value = 58.142857
$ws.Range("K7").Value = 58.142857
$ws.Range("M7").Value = 54.857143
$ws.Range("H16").Value = 1655.8334
$ws.Range("I16").Value = 1695
$ws.Range("J16").Value = 1616.6666
$ws.Range("K16").Value = 1695
$ws.Range("L16").Value = 1616.6666
$ws.Range("M16").Value = -1408
$ws.Range("N16").Value = -2190.6666
$ws.Range("H36").Value = 15000
$ws.Range("I36").Value = 15000
$ws.Range("K36").Value = 15000
$ws.Range("M36").Value = -14612
$ws.Range("H39").Value = 6017
$ws.Range("I39").Value = 4025.5
$ws.Range("K39").Value = 4025.5
$ws.Range("M39").Value = -3634.5
$ws.Range("H40").Value = 15000
$ws.Range("I40").Value = 15000
$ws.Range("K40").Value = 15000
$ws.Range("M40").Value = -14840
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9572
$ws.Range("H49").Value = 6017
$ws.Range("I49").Value = 4025.5
$ws.Range("K49").Value = 4025.5
$ws.Range("M49").Value = -3843.5
$ws.Range("H50").Value = 19333.334
$ws.Range("H51").Value = 40000
$ws.Range("H59").Value = 38750
$ws.Range("I59").Value = 50000
$ws.Range("K59").Value = 50000
$ws.Range("M59").Value = -48855
$ws.Range("H61").Value = 40000
$ws.Range("H113").Value = 1655.8334
$ws.Range("I113").Value = 1695
$ws.Range("J113").Value = 1616.6666
$ws.Range("K113").Value = 1695
$ws.Range("L113").Value = 1616.6666
$ws.Range("M113").Value = 475
$ws.Range("N113").Value = -5956.6666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 756.3333
$ws.Range("I103").Value = 1439.6666
$ws.Range("J103").Value = 414.66666
$ws.Range("K103").Value = 4318.9998
$ws.Range("L103").Value = 1243.99998
$ws.Range("M103").Value = -3439.9998
$ws.Range("N103").Value = -3001.99998
$ws.Range("H141").Value = 4140
$ws.Range("I141").Value = 3190
$ws.Range("J141").Value = 4900
$ws.Range("K141").Value = 9570
$ws.Range("L141").Value = 14700
$ws.Range("M141").Value = -4390
$ws.Range("N141").Value = -25060

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4395.8687
$ws.Range("I102").Value = 5073.6206
$ws.Range("K102").Value = 5073.6206
$ws.Range("M102").Value = -3451.6206
$ws.Range("H126").Value = 13319.518
$ws.Range("I126").Value = 22701.25
$ws.Range("J126").Value = 9745.522999999999
$ws.Range("K126").Value = 68103.75
$ws.Range("L126").Value = 29236.569
$ws.Range("M126").Value = -65633.75
$ws.Range("N126").Value = -34176.569
$ws.Range("H132").Value = 3012.8235
$ws.Range("I132").Value = 2381
$ws.Range("J132").Value = 4033.4614
$ws.Range("K132").Value = 7143
$ws.Range("L132").Value = 12100.3842
$ws.Range("M132").Value = -4613
$ws.Range("N132").Value = -17160.3842
$ws.Range("H134").Value = 40229.8
$ws.Range("J134").Value = 40229.8
$ws.Range("L134").Value = 120689.4
$ws.Range("N134").Value = -125759.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 20834.408
$ws.Range("I40").Value = 21849.521
$ws.Range("K40").Value = 21849.521
$ws.Range("M40").Value = -21713.521
$ws.Range("H61").Value = 599
$ws.Range("I61").Value = 599
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 599
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -397
$ws.Range("H113").Value = 599
$ws.Range("I113").Value = 599
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 599
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1571
$ws.Range("H136").Value = 7582.5386
$ws.Range("I136").Value = 3368.8572
$ws.Range("J136").Value = 12498.5
$ws.Range("K136").Value = 10106.5716
$ws.Range("L136").Value = 37495.5
$ws.Range("M136").Value = -7556.571599999999
$ws.Range("N136").Value = -42595.5
$ws.Range("N61").ClearContents()
$ws.Range("N113").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 13336
$ws.Range("J30").Value = 9999.5
$ws.Range("L30").Value = 9999.5
$ws.Range("N30").Value = -10213.5
$ws.Range("H107").Value = 37049.223
$ws.Range("I107").Value = 3920
$ws.Range("K107").Value = 11760
$ws.Range("M107").Value = -9840
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("H132").Value = 7266.5117
$ws.Range("I132").Value = 10359.675
$ws.Range("J132").Value = 4454.5454
$ws.Range("K132").Value = 31079.025
$ws.Range("L132").Value = 13363.6362
$ws.Range("M132").Value = -28549.025
$ws.Range("M131").ClearContents()
